$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated list of names for rows 2-23 (column B), column A holds the
# sequential index (1-based) which is just (row-1).
$names = @(
    "AANYA JAIN",
    "AARAV DUA",
    "AAYUSH GUPTA",
    "AKSHITA PURI",
    "ANSHIKA",
    "ARNAV SHARMA",
    "AROUSH SETH",
    "BHAVYA SHARMA",
    "DIPIN PANDEY",
    "ISHANI JHA",
    "LAKSHAY MALHOTRA",
    "KASHIKA TAYAL",
    "MAHI WADHWA",
    "PANKAJ",
    "PRATHAM SHARMA",
    "RANVEER SOLANKI",
    "RISHABH SINGH",
    "ROUNAK BISWAS",
    "SUMAN",
    "UNNABH BHALLA",
    "VANSHIKA ARYA",
    "YUVRAJ MALIK"
)

# First, extend column A formatting (border/bold/center/top alignment)
# down to the new rows 17-23 by copying the format of the last existing
# formatted cell (A16).
$ws.Range("A16").Copy() | Out-Null
$ws.Range("A17:A23").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# Now write the sequence numbers (1..22) into A2:A23 and the names into
# B2:B23.
for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $i + 1
    $ws.Cells.Item($row, 2).Value = $names[$i]
}

Write-Host "Updated attendance list through row $($names.Length + 1)"
